$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ELC_TRADE")

# Change prices of Bornholm (DKISLBH) and ISL1 (DKISL1) trades.
# Row 11 = INVCOST for DKISLBH, Row 12 = FIXOM for DKISLBH (col H, mirrored into col J)
$ws.Range("H11").Value = 0.1
$ws.Range("H12").Value = 0.1

# Row 17 = INVCOST for DKISL1, Row 18 = FIXOM for DKISL1 (col I, mirrored into col K)
$ws.Range("I17").Value = 0.1
$ws.Range("I18").Value = 0.1

# Re-enter the mirroring formulas so Excel regroups them into shared formulas
# (J10:J15 mirrors H10:H15, K16:K21 mirrors I16:I21), same as the authored workbook.
$ws.Range("J10:J15").Formula = "=H10"
$ws.Range("K16:K21").Formula = "=I16"

# Update the active selection to match the diff
$ws.Range("K21").Select()
